$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sessions")

# Collapse the per-clip start/end column pairs into single clip_out_N /
# clip_in_N columns (clip_out_start_N + clip_out_end_N -> clip_out_N,
# clip_in_start_N + clip_in_end_N -> clip_in_N) for both clip 1 and clip 2,
# shifting the surviving headers (file_2, clip_out_2, clip_in_2) left.
$ws.Range("T1").Value = "clip_out_1"
$ws.Range("U1").Value = "clip_in_1"
$ws.Range("V1").Value = "file_2"
$ws.Range("W1").Value = "clip_out_2"
$ws.Range("X1").Value = "clip_in_2"

# The old Y1:AB1 headers (clip_out_start_2, clip_out_end_2,
# clip_in_start_2, clip_in_end_2) no longer have a place - drop them so the
# sheet's used range shrinks back down to A1:X1.
$ws.Range("Y1:AB1").ClearContents()
